$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3365.4  # was 3377.1143
$ws.Range("I64").Value = 3142.7144  # was 3112.375
$ws.Range("J64").Value = 3421.0715  # was 3455.5557
$ws.Range("K64").Value = 3142.7144  # was 3112.375
$ws.Range("L64").Value = 3421.0715  # was 3455.5557
$ws.Range("M64").Value = -2894.7144  # was -2864.375
$ws.Range("N64").Value = -3917.0715  # was -3951.5557
$ws.Range("H67").Value = 3365.4  # was 3377.1143
$ws.Range("I67").Value = 3142.7144  # was 3112.375
$ws.Range("J67").Value = 3421.0715  # was 3455.5557
$ws.Range("K67").Value = 3142.7144  # was 3112.375
$ws.Range("L67").Value = 3421.0715  # was 3455.5557
$ws.Range("M67").Value = -2284.7144  # was -2254.375
$ws.Range("N67").Value = -5137.0715  # was -5171.5557
$ws.Range("H70").Value = 63893.938  # was 84750
$ws.Range("I70").Value = 501150  # was 334433.34
$ws.Range("J70").Value = 1428.7858  # was 1522.2222
$ws.Range("K70").Value = 1503450  # was 1003300.02
$ws.Range("L70").Value = 4286.357400000001  # was 4566.6666
$ws.Range("M70").Value = -1503180  # was -1003030.02
$ws.Range("N70").Value = -4826.357400000001  # was -5106.6666
$ws.Range("H73").Value = 63893.938  # was 84750
$ws.Range("I73").Value = 501150  # was 334433.34
$ws.Range("J73").Value = 1428.7858  # was 1522.2222
$ws.Range("K73").Value = 1503450  # was 1003300.02
$ws.Range("L73").Value = 4286.357400000001  # was 4566.6666
$ws.Range("M73").Value = -1502514  # was -1002364.02
$ws.Range("N73").Value = -6158.357400000001  # was -6438.6666
$ws.Range("H74").Value = 4379  # was 4220.647
$ws.Range("I74").Value = 3975.75  # was 3841
$ws.Range("K74").Value = 3975.75  # was 3841
$ws.Range("M74").Value = -3039.75  # was -2905
$ws.Range("H77").Value = 4379  # was 4220.647
$ws.Range("I77").Value = 3975.75  # was 3841
$ws.Range("K77").Value = 19878.75  # was 19205
$ws.Range("M77").Value = -15198.75  # was -14525
$ws.Range("H101").Value = 534.9231  # was 2099.6428
$ws.Range("I101").Value = 556  # was 2451
$ws.Range("J101").Value = 487.5  # was 1221.25
$ws.Range("K101").Value = 1668  # was 7353
$ws.Range("L101").Value = 1462.5  # was 3663.75
$ws.Range("M101").Value = -46  # was -5731
$ws.Range("N101").Value = -4706.5  # was -6907.75
$ws.Range("H111").Value = 1879.2222  # was 2121
$ws.Range("I111").Value = 1755.2941  # was 2128.9285
$ws.Range("J111").Value = 2089.9  # was 2109.9
$ws.Range("K111").Value = 5265.8823  # was 6386.7855
$ws.Range("L111").Value = 6269.700000000001  # was 6329.700000000001
$ws.Range("M111").Value = -2198.8823  # was -3319.7855
$ws.Range("N111").Value = -12403.7  # was -12463.7
$ws.Range("H112").Value = 2425.7532  # was 2334.092
$ws.Range("I112").Value = 325  # was 300
$ws.Range("J112").Value = 2534.883  # was 2432.1204
$ws.Range("K112").Value = 975  # was 900
$ws.Range("L112").Value = 7604.648999999999  # was 7296.361199999999
$ws.Range("M112").Value = 133  # was 208
$ws.Range("N112").Value = -9820.648999999999  # was -9512.361199999999
$ws.Range("H113").Value = 2957.2144  # was 3001.2307
$ws.Range("I113").Value = 2413.5715  # was 2435
$ws.Range("J113").Value = 3500.8572  # was 3486.5715
$ws.Range("K113").Value = 2413.5715  # was 2435
$ws.Range("L113").Value = 3500.8572  # was 3486.5715
$ws.Range("M113").Value = 840.4285  # was 819
$ws.Range("N113").Value = -10008.8572  # was -9994.5715
$ws.Range("H115").Value = 1278.1428  # was 2000
$ws.Range("I115").Value = 676.4  # was 2000
$ws.Range("J115").Value = 2782.5  # was 0
$ws.Range("K115").Value = 2029.2  # was 6000
$ws.Range("L115").Value = 8347.5  # was 0
$ws.Range("M115").Value = -462.1999999999998  # was -4433
$ws.Range("N115").Value = -11481.5  # was None
$ws.Range("H129").Value = 1191.5667  # was 1222.5769
$ws.Range("J129").Value = 1267.6923  # was 1318.1818
$ws.Range("L129").Value = 3803.0769  # was 3954.5454
$ws.Range("N129").Value = -13803.0769  # was -13954.5454
$ws.Range("H138").Value = 2655.738  # was 2521.0254
$ws.Range("I138").Value = 1573.7609  # was 1487.44
$ws.Range("J138").Value = 3965.5  # was 4303.069
$ws.Range("K138").Value = 4721.2827  # was 4462.32
$ws.Range("L138").Value = 11896.5  # was 12909.207
$ws.Range("M138").Value = 418.7173000000003  # was 677.6800000000003
$ws.Range("N138").Value = -22176.5  # was -23189.207

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0  # was 1197
$ws.Range("I4").Value = 0  # was 1197
$ws.Range("K4").Value = 0  # was 1197
$ws.Range("M4").ClearContents()  # was -1081
$ws.Range("H105").Value = 275185  # was 500370
$ws.Range("J105").Value = 275185  # was 500370
$ws.Range("L105").Value = 275185  # was 500370
$ws.Range("N105").Value = -282173  # was -507358
$ws.Range("H132").Value = 1802.5  # was 1871.4714
$ws.Range("I132").Value = 1423.8518  # was 1490.12
$ws.Range("K132").Value = 4271.555399999999  # was 4470.36
$ws.Range("M132").Value = -1741.555399999999  # was -1940.36

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 927568.9  # was 843292.75
$ws.Range("I58").Value = 1482998.2  # was 1373184.9
$ws.Range("J58").Value = 1853.2667  # was 1699.2354
$ws.Range("K58").Value = 1482998.2  # was 1373184.9
$ws.Range("L58").Value = 1853.2667  # was 1699.2354
$ws.Range("M58").Value = -1482795.2  # was -1372981.9
$ws.Range("N58").Value = -2259.2667  # was -2105.2354
$ws.Range("H105").Value = 7126.25  # was 8801.462
$ws.Range("I105").Value = 9792.727999999999  # was 10742
$ws.Range("J105").Value = 1260  # was 2333
$ws.Range("K105").Value = 9792.727999999999  # was 10742
$ws.Range("L105").Value = 1260  # was 2333
$ws.Range("M105").Value = -8045.727999999999  # was -8995
$ws.Range("N105").Value = -4754  # was -5827
$ws.Range("H136").Value = 927568.9  # was 843292.75
$ws.Range("I136").Value = 1482998.2  # was 1373184.9
$ws.Range("J136").Value = 1853.2667  # was 1699.2354
$ws.Range("K136").Value = 4448994.6  # was 4119554.7
$ws.Range("L136").Value = 5559.800099999999  # was 5097.706200000001
$ws.Range("M136").Value = -4446444.6  # was -4117004.7
$ws.Range("N136").Value = -10659.8001  # was -10197.7062
$ws.Range("H139").Value = 20000  # was 23500
$ws.Range("J139").Value = 0  # was 24666.666
$ws.Range("L139").Value = 0  # was 24666.666
$ws.Range("N139").ClearContents()  # was -34946.666

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 666372.8  # was 623381.9399999999
$ws.Range("I12").Value = 56.625  # was 51.333332
$ws.Range("J12").Value = 920207.5600000001  # was 878380.9
$ws.Range("K12").Value = 169.875  # was 153.999996
$ws.Range("L12").Value = 2760622.68  # was 2635142.7
$ws.Range("M12").Value = 3.125  # was 19.00000399999999
$ws.Range("N12").Value = -2760968.68  # was -2635488.7
$ws.Range("H82").Value = 13946.667  # was 12902
$ws.Range("I82").Value = 0  # was 1000
$ws.Range("J82").Value = 13946.667  # was 14224.444
$ws.Range("K82").Value = 0  # was 3000
$ws.Range("L82").Value = 41840.001  # was 42673.33199999999
$ws.Range("M82").ClearContents()  # was -2594
$ws.Range("N82").Value = -42652.001  # was -43485.33199999999
$ws.Range("H85").Value = 13946.667  # was 12902
$ws.Range("I85").Value = 0  # was 1000
$ws.Range("J85").Value = 13946.667  # was 14224.444
$ws.Range("K85").Value = 0  # was 3000
$ws.Range("L85").Value = 41840.001  # was 42673.33199999999
$ws.Range("M85").ClearContents()  # was -1596
$ws.Range("N85").Value = -44648.001  # was -45481.33199999999
$ws.Range("H122").Value = 850  # was 745.6667
$ws.Range("I122").Value = 487.25  # was 497.25
$ws.Range("J122").Value = 1264.5714  # was 869.875
$ws.Range("K122").Value = 4385.25  # was 4475.25
$ws.Range("L122").Value = 11381.1426  # was 7828.875
$ws.Range("M122").Value = -1935.25  # was -2025.25
$ws.Range("N122").Value = -16281.1426  # was -12728.875
$ws.Range("H131").Value = 962.89  # was 984.1799999999999
$ws.Range("I131").Value = 375.7143  # was 388.33334
$ws.Range("J131").Value = 1007.086  # was 1022.21277
$ws.Range("K131").Value = 1127.1429  # was 1165.00002
$ws.Range("L131").Value = 3021.258  # was 3066.63831
$ws.Range("M131").Value = 3912.8571  # was 3874.99998
$ws.Range("N131").Value = -13101.258  # was -13146.63831
$ws.Range("H133").Value = 4490.0586  # was 4691.1665
$ws.Range("I133").Value = 2431.8333  # was 2675.25
$ws.Range("J133").Value = 5612.727  # was 5267.143
$ws.Range("K133").Value = 7295.499899999999  # was 8025.75
$ws.Range("L133").Value = 16838.181  # was 15801.429
$ws.Range("M133").Value = -2235.499899999999  # was -2965.75
$ws.Range("N133").Value = -26958.181  # was -25921.429
$ws.Range("H134").Value = 5255.5  # was 5171.926
$ws.Range("J134").Value = 6031.8335  # was 5872.2104
$ws.Range("L134").Value = 18095.5005  # was 17616.6312
$ws.Range("N134").Value = -28235.5005  # was -27756.6312
$ws.Range("H137").Value = 3542.5715  # was 2564.8
$ws.Range("I137").Value = 2000  # was 1623
$ws.Range("J137").Value = 3799.6667  # was 3977.5
$ws.Range("K137").Value = 6000  # was 4869
$ws.Range("L137").Value = 11399.0001  # was 11932.5
$ws.Range("M137").Value = -900  # was 231
$ws.Range("N137").Value = -21599.0001  # was -22132.5
$ws.Range("H138").Value = 1042.25  # was 2224.9092
$ws.Range("I138").Value = 1042.25  # was 993.6
$ws.Range("J138").Value = 0  # was 3251
$ws.Range("K138").Value = 3126.75  # was 2980.8
$ws.Range("L138").Value = 0  # was 9753
$ws.Range("M138").Value = 2013.25  # was 2159.2
$ws.Range("N138").ClearContents()  # was -20033
$ws.Range("H139").Value = 2927.0908  # was 2175.9714
$ws.Range("I139").Value = 2881.111  # was 1849.56
$ws.Range("J139").Value = 2958.923  # was 2992
$ws.Range("K139").Value = 8643.332999999999  # was 5548.68
$ws.Range("L139").Value = 8876.769  # was 8976
$ws.Range("M139").Value = -3503.332999999999  # was -408.6800000000003
$ws.Range("N139").Value = -19156.769  # was -19256
$ws.Range("H140").Value = 5826.65  # was 1483.8334
$ws.Range("I140").Value = 1150  # was 1199.0869
$ws.Range("J140").Value = 9653  # was 8033
$ws.Range("K140").Value = 3450  # was 3597.2607
$ws.Range("L140").Value = 28959  # was 24099
$ws.Range("M140").Value = 1730  # was 1582.7393
$ws.Range("N140").Value = -39319  # was -34459

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1066.6666  # was 924
$ws.Range("I22").Value = 850  # was 706.6667
$ws.Range("J22").Value = 1500  # was 1250
$ws.Range("K22").Value = 850  # was 706.6667
$ws.Range("L22").Value = 1500  # was 1250
$ws.Range("M22").Value = -555  # was -411.6667
$ws.Range("N22").Value = -2090  # was -1840
$ws.Range("H27").Value = 1066.6666  # was 924
$ws.Range("I27").Value = 850  # was 706.6667
$ws.Range("J27").Value = 1500  # was 1250
$ws.Range("K27").Value = 850  # was 706.6667
$ws.Range("L27").Value = 1500  # was 1250
$ws.Range("M27").Value = -743  # was -599.6667
$ws.Range("N27").Value = -1714  # was -1464
$ws.Range("H136").Value = 1556.4559  # was 1495.7595
$ws.Range("I136").Value = 1322.9678  # was 1260.7858
$ws.Range("J136").Value = 3969.1667  # was 3323.3333
$ws.Range("K136").Value = 3968.9034  # was 3782.3574
$ws.Range("L136").Value = 11907.5001  # was 9969.999899999999
$ws.Range("M136").Value = -1418.9034  # was -1232.3574
$ws.Range("N136").Value = -17007.5001  # was -15069.9999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 50000  # was 50002
$ws.Range("J18").Value = 0  # was 50002.332
$ws.Range("L18").Value = 0  # was 50002.332
$ws.Range("N18").ClearContents()  # was -50348.332
$ws.Range("H113").Value = 508.25  # was 598.875
$ws.Range("I113").Value = 443.14285  # was 483.5
$ws.Range("J113").Value = 660.1667  # was 945
$ws.Range("K113").Value = 1329.42855  # was 1450.5
$ws.Range("L113").Value = 1980.5001  # was 2835
$ws.Range("M113").Value = 840.5714499999999  # was 719.5
$ws.Range("N113").Value = -6320.5001  # was -7175
$ws.Range("H132").Value = 1264.5869  # was 1818.875
$ws.Range("I132").Value = 863.82355  # was 1590.4706
$ws.Range("J132").Value = 2400.0833  # was 2373.5715
$ws.Range("K132").Value = 2591.47065  # was 4771.4118
$ws.Range("L132").Value = 7200.249899999999  # was 7120.7145
$ws.Range("M132").Value = -61.47064999999975  # was -2241.4118
$ws.Range("N132").Value = -12260.2499  # was -12180.7145
